# ---------------------------------------------------------------------------
# C5-PowerPoint.pptx edit
#
# 1) The table on slide 6 gets switched from the deck's custom "Table_0"
#    style to the built-in table style {136B0D21-0B51-4B61-8D5B-5710C4329CC6}.
#
# 2) The design/theme applied to the deck is swapped from "Integral" to the
#    default "Office Theme" palette (the deck's notesMaster keeps carrying
#    the old, now-unused "Office Theme" theme part, while the live theme
#    used by the slide master takes on the Office Theme color values).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Retarget the table's style on slide 6 --------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{136B0D21-0B51-4B61-8D5B-5710C4329CC6}")
    }
}

# --- 2) Swap the active theme's colors from "Integral" to "Office Theme" ----
function HexToRgbInt([string]$hex) {
    $v = [Convert]::ToInt32($hex, 16)
    $r = ($v -shr 16) -band 0xFF
    $g = ($v -shr 8) -band 0xFF
    $b = $v -band 0xFF
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

# Office Theme color scheme, in ThemeColorScheme.Item(1..12) slot order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = HexToRgbInt($officeThemeColors[$i - 1])
}
